$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find a phrase in the document body and return its Range (collapsed
# to the match). Throws if not found, so mistakes fail loudly.
# ---------------------------------------------------------------------------
function Find-Range([string]$text) {
    $r = $d.Content
    $found = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $text"
    }
    return $r
}

# ---------------------------------------------------------------------------
# Edit 2 (done first so it does not disturb whitespace-handling of the run
# split done by Edit 1 below): move the "_GoBack" last-edit bookmark to just
# after "Additionally the player " (i.e. right before "will be able to
# interact"). Word only allows one bookmark per name, so adding a new
# "_GoBack" moves it off of its old location automatically.
# ---------------------------------------------------------------------------
$rAfterPlayer = Find-Range("Additionally the player ")
$goBackPos = $rAfterPlayer.End
$rGoBack = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $rGoBack)

# ---------------------------------------------------------------------------
# Edit 1: "...who will spawn consistent waves..." -> "...spawn constant waves..."
# Word replaces the single word; to keep it as its own run (as the real
# edit session produced) we touch formatting on the replacement so it stays
# split from its neighbours instead of being silently re-merged.
# ---------------------------------------------------------------------------
$rConsistent = Find-Range("consistent")
$rConsistent.Text = "constant"
$start = $rConsistent.Start
$rNewWord = $d.Range($start, $start + 8)   # "constant" is 8 chars
$rNewWord.Bold = 1
$rNewWord.Bold = 0

# ---------------------------------------------------------------------------
# Edit 3: "Will also take lead in creating the colliders, meshes, and
# scripts for the environment..." -> "...the scripted events for the
# environment..."
# ---------------------------------------------------------------------------
$rColliders1 = Find-Range("also take lead in creating the colliders, meshes, and scripts")
$rColliders1.Text = "also take lead in creating the scripted events"
$s1 = $rColliders1.Start
$rNew1 = $d.Range($s1, $s1 + ("also take lead in creating the scripted events").Length)
$rNew1.Bold = 1
$rNew1.Bold = 0

# ---------------------------------------------------------------------------
# Edit 4: "Character design lead is responsible for creating colliders,
# meshes, and scripts for both player..." -> "...creating the scripted
# events for both player..."
# ---------------------------------------------------------------------------
$rColliders2 = Find-Range("is responsible for creating colliders, meshes, and scripts")
$rColliders2.Text = "is responsible for creating the scripted events"
$s2 = $rColliders2.Start
# Only "the scripted events" is its own run in the target; recompute its
# start based on the fixed prefix length.
$prefix2 = "is responsible for creating "
$s2b = $s2 + $prefix2.Length
$rNew2 = $d.Range($s2b, $s2b + ("the scripted events").Length)
$rNew2.Bold = 1
$rNew2.Bold = 0

# ---------------------------------------------------------------------------
# Edit 5: remove the old "_GoBack" bookmark split in the last paragraph,
# re-merging "...other f" + "iles." into a single run "...other files."
# Running Find/Replace (rather than re-assigning .Text) over the unchanged
# sentence collapses the two runs and drops the now-stale bookmark markers
# (Edit 2 above already re-homed the "_GoBack" bookmark name itself).
# ---------------------------------------------------------------------------
$rFilesScope = $d.Content
$foundFiles = $rFilesScope.Find.Execute( `
    "Free open source version control system to share source code, assets and other files.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Free open source version control system to share source code, assets and other files.", 2)
Write-Output "files paragraph normalized: $foundFiles"

# ---------------------------------------------------------------------------
# Edit 6: header page-number field cached result "5" -> "4".
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$rHdr = $hdr.Range
$foundPage = $rHdr.Find.Execute("5", $true, $false, $false, $false, $false, $true, 1, $false, "4", 2)
Write-Output "page field updated: $foundPage"

Write-Output "done"
